$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.231550931930542
$ws.Range("B1").Value = 2.154545783996582
$ws.Range("C1").Value = 6.016922950744629
$ws.Range("D1").Value = 1.980476379394531
$ws.Range("E1").Value = 1.150825262069702
